$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new numeric-looking values would
# otherwise be auto-converted to numbers by Excel, so they stay text
# (matching the original inlineStr cell type).
$textCells = @("D14", "D38", "D39", "D9", "D42", "D5", "D30", "D10", "D48", "D23", "D36", "D17", "D21", "D50", "D49", "D27", "D44", "D16", "D37", "D25", "D13", "D32", "D34", "D8", "D28", "D11", "D26", "D19", "D35", "D47", "D6")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '28.623.76'
$ws.Range("E2").Value = '  -2.15%  '
$ws.Range("D3").Value = '1.796.42'
$ws.Range("E3").Value = '  -1.84%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '232.27'
$ws.Range("E5").Value = '  -1.15%  '
$ws.Range("D6").Value = '0.5903'
$ws.Range("E6").Value = '  -1.78%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.2779'
$ws.Range("E8").Value = '  +0.30%  '
$ws.Range("D9").Value = '23.36'
$ws.Range("E9").Value = '  -0.32%  '
$ws.Range("D10").Value = '0.06764'
$ws.Range("E10").Value = '  -3.21%  '
$ws.Range("D11").Value = '0.07560'
$ws.Range("E11").Value = '  -1.05%  '
$ws.Range("D12").Value = '1.799.83'
$ws.Range("E12").Value = '  -1.70%  '
$ws.Range("D13").Value = '4.798'
$ws.Range("E13").Value = '  +0.63%  '
$ws.Range("D14").Value = '0.6149'
$ws.Range("E14").Value = '  -1.76%  '
$ws.Range("D15").Value = '2.039.51'
$ws.Range("E15").Value = '  -1.90%  '
$ws.Range("D16").Value = '75.61'
$ws.Range("E16").Value = '  -3.89%  '
$ws.Range("D17").Value = '0.000008922'
$ws.Range("E17").Value = '  -8.16%  '
$ws.Range("D18").Value = '28.612.69'
$ws.Range("E18").Value = '  -1.70%  '
$ws.Range("D19").Value = '5.411'
$ws.Range("E19").Value = '  -6.09%  '
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").Value = '209.37'
$ws.Range("E21").Value = '  -5.87%  '
$ws.Range("E22").Value = '  -0.90%  '
$ws.Range("D23").Value = '6.839'
$ws.Range("E23").Value = '  -1.10%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").Value = '152.52'
$ws.Range("E25").Value = '  -2.44%  '
$ws.Range("D26").Value = '8.120'
$ws.Range("E26").Value = '  +1.91%  '
$ws.Range("D27").Value = '0.1264'
$ws.Range("E27").Value = '  -2.82%  '
$ws.Range("D28").Value = '16.45'
$ws.Range("E28").Value = '  -0.66%  '
$ws.Range("E29").Value = '  -2.71%  '
$ws.Range("D30").Value = '0.06220'
$ws.Range("E30").Value = '  -10.35%  '
$ws.Range("E31").Value = '  -1.38%  '
$ws.Range("D32").Value = '3.797'
$ws.Range("E32").Value = '  -0.98%  '
$ws.Range("E33").Value = '  +0.53%  '
$ws.Range("D34").Value = '1.738'
$ws.Range("E34").Value = '  +0.85%  '
$ws.Range("D35").Value = '1.049'
$ws.Range("E35").Value = '  -4.47%  '
$ws.Range("D36").Value = '0.6431'
$ws.Range("E36").Value = '  -0.30%  '
$ws.Range("D37").Value = '2.503'
$ws.Range("E37").Value = '  -1.58%  '
$ws.Range("D38").Value = '2.720'
$ws.Range("E38").Value = '  -0.79%  '
$ws.Range("D39").Value = '0.01700'
$ws.Range("E39").Value = '  -2.46%  '
$ws.Range("E40").Value = '  -2.33%  '
$ws.Range("D41").Value = '1.148.88'
$ws.Range("E41").Value = '  -4.52%  '
$ws.Range("D42").Value = '0.8768'
$ws.Range("E42").Value = '  -3.04%  '
$ws.Range("E43").Value = '  +0.09%  '
$ws.Range("D44").Value = '100.38'
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("D45").Value = '1.947.85'
$ws.Range("E45").Value = '  -1.91%  '
$ws.Range("E46").Value = '  -3.25%  '
$ws.Range("D47").Value = '0.00000000111'
$ws.Range("E47").Value = '  -2.88%  '
$ws.Range("D48").Value = '1.588'
$ws.Range("E48").Value = '  +0.63%  '
$ws.Range("D49").Value = '8.353'
$ws.Range("E49").Value = '  -1.82%  '
$ws.Range("D50").Value = '0.05461'
$ws.Range("E50").Value = '  -0.84%  '
$ws.Range("E51").Value = '  -1.67%  '
